# Update the "Last Updated: 21-Apr-20" date stamp to "23-Apr-20" on every
# slide of the deck. The date appears as the final run inside the
# "TextBox 2" shape's single paragraph (after "Project name: ...",
# "Project manager: ...", "Last Updated: "), so we locate the substring
# within the full TextRange and replace just that run's characters —
# this keeps every other run (and its formatting/hyperlink) untouched.

$p = $ppt.ActivePresentation

$oldDate = "21-Apr-20"
$newDate = "23-Apr-20"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $fullText = $tr.Text

            if ($fullText.Contains($oldDate)) {
                $idx = $fullText.IndexOf($oldDate)
                $run = $tr.Characters($idx + 1, $oldDate.Length)
                $run.Text = $newDate
            }
        }
    }
}
